# Update the "想去人数" (want-to-go count) figures in column F for the
# rows that changed in the source data refresh (commit: "Update gh-pages
# to output generated at 456a3b4"). The same rows are updated identically
# on both the "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

$updates = @{
    5  = 251
    6  = 41
    15 = 435
    17 = 477
    18 = 403
    23 = 1107
    24 = 2835
    28 = 49
    29 = 1609
    31 = 453
    35 = 602
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
